$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1823923333333334
$ws.Range("H2").Value = 0.547177
$ws.Range("M2").Value = 0.4214143333333333
$ws.Range("N2").Value = 1.264243
$ws.Range("O2").Value = 0.02434128610922473
$ws.Range("P2").Value = 0.02434128610922473
$ws.Range("Q2").Value = 0.07686274355677779
$ws.Range("R2").Value = 0.691764692011
$ws.Range("S2").Value = 0.02434128610922473
$ws.Range("T2").Value = 0.02434128610922473

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1823923333333334
$ws.Range("H3").Value = 0.547177
$ws.Range("N3").Value = 36.386704
$ws.Range("O3").Value = 0.7005766871049885
$ws.Range("P3").Value = 0.7005766871049887
$ws.Range("Q3").Value = 2.212218614956444
$ws.Range("R3").Value = 19.909967534608
$ws.Range("S3").Value = 0.7005766871049885
$ws.Range("T3").Value = 0.7005766871049887

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1823923333333334
$ws.Range("H4").Value = 0.547177
$ws.Range("M4").Value = 4.762423333333333
$ws.Range("N4").Value = 14.28727
$ws.Range("O4").Value = 0.2750820267857866
$ws.Range("P4").Value = 0.2750820267857866
$ws.Range("Q4").Value = 0.8686295040877778
$ws.Range("R4").Value = 7.81766553679
$ws.Range("S4").Value = 0.2750820267857866
$ws.Range("T4").Value = 0.2750820267857866
